# Auto-generated edit script applying the Pandaemonium_Profits diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3527
$ws.Range("J69").Value = 3527
$ws.Range("L69").Value = 10581
$ws.Range("N69").Value = -12329

$ws.Range("H72").Value = 3527
$ws.Range("J72").Value = 3527
$ws.Range("L72").Value = 31743
$ws.Range("N72").Value = -40479

$ws.Range("H98").Value = 717.5
$ws.Range("I98").Value = 511.5625
$ws.Range("J98").Value = 1266.6666
$ws.Range("K98").Value = 511.5625
$ws.Range("L98").Value = 1266.6666
$ws.Range("M98").Value = 986.4375
$ws.Range("N98").Value = -4262.6666

$ws.Range("H106").Value = 2598.6316
$ws.Range("I106").Value = 2578.9092
$ws.Range("K106").Value = 2578.9092
$ws.Range("M106").Value = -1947.9092

$ws.Range("H116").Value = 2052
$ws.Range("I116").Value = 1795.9231
$ws.Range("J116").Value = 2421.889
$ws.Range("K116").Value = 1795.9231
$ws.Range("L116").Value = 2421.889
$ws.Range("M116").Value = 1646.0769
$ws.Range("N116").Value = -9305.888999999999

$ws.Range("H122").Value = 717.5
$ws.Range("I122").Value = 511.5625
$ws.Range("J122").Value = 1266.6666
$ws.Range("K122").Value = 1534.6875
$ws.Range("L122").Value = 3799.9998
$ws.Range("M122").Value = 915.3125
$ws.Range("N122").Value = -8699.9998

$ws.Range("H125").Value = 778.8946999999999
$ws.Range("J125").Value = 799.94116
$ws.Range("L125").Value = 7199.47044
$ws.Range("N125").Value = -12119.47044

$ws.Range("H129").Value = 1025.5231
$ws.Range("J129").Value = 1053.9517
$ws.Range("L129").Value = 3161.8551
$ws.Range("N129").Value = -13161.8551

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28149.592
$ws.Range("I32").Value = 32919.266
$ws.Range("K32").Value = 32919.266
$ws.Range("M32").Value = -32632.266

$ws.Range("H61").Value = 5288.7
$ws.Range("I61").Value = 4120.9614
$ws.Range("J61").Value = 7457.357
$ws.Range("K61").Value = 4120.9614
$ws.Range("L61").Value = 7457.357
$ws.Range("M61").Value = -3908.9614
$ws.Range("N61").Value = -7881.357

$ws.Range("H110").Value = 1394.7
$ws.Range("I110").Value = 1407.4667
$ws.Range("J110").Value = 1356.4
$ws.Range("K110").Value = 1407.4667
$ws.Range("L110").Value = 1356.4
$ws.Range("M110").Value = 637.5333000000001
$ws.Range("N110").Value = -5446.4

$ws.Range("H132").Value = 2472.8367
$ws.Range("I132").Value = 2050.182
$ws.Range("J132").Value = 2817.2222
$ws.Range("K132").Value = 6150.545999999999
$ws.Range("L132").Value = 8451.6666
$ws.Range("M132").Value = -3620.545999999999
$ws.Range("N132").Value = -13511.6666

$ws.Range("H136").Value = 5288.7
$ws.Range("I136").Value = 4120.9614
$ws.Range("J136").Value = 7457.357
$ws.Range("K136").Value = 12362.8842
$ws.Range("L136").Value = 22372.071
$ws.Range("M136").Value = -9812.8842
$ws.Range("N136").Value = -27472.071

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1286.45
$ws.Range("I94").Value = 1132.4375
$ws.Range("J94").Value = 1902.5
$ws.Range("K94").Value = 1132.4375
$ws.Range("L94").Value = 1902.5
$ws.Range("M94").Value = -681.4375
$ws.Range("N94").Value = -2804.5

$ws.Range("H135").Value = 38833.332
$ws.Range("J135").Value = 38833.332
$ws.Range("L135").Value = 38833.332
$ws.Range("N135").Value = -48973.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2757629.5
$ws.Range("I58").Value = 6062984
$ws.Range("J58").Value = 3167.2222
$ws.Range("K58").Value = 6062984
$ws.Range("L58").Value = 3167.2222
$ws.Range("M58").Value = -6062781
$ws.Range("N58").Value = -3573.2222

$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2376
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11880
$ws.Range("N65").ClearContents()

$ws.Range("H136").Value = 2757629.5
$ws.Range("I136").Value = 6062984
$ws.Range("J136").Value = 3167.2222
$ws.Range("K136").Value = 18188952
$ws.Range("L136").Value = 9501.6666
$ws.Range("M136").Value = -18186402
$ws.Range("N136").Value = -14601.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6466.5835
$ws.Range("I87").Value = 2499.6667
$ws.Range("J87").Value = 7788.8887
$ws.Range("K87").Value = 7499.000100000001
$ws.Range("L87").Value = 23366.6661
$ws.Range("M87").Value = -6251.000100000001
$ws.Range("N87").Value = -25862.6661

$ws.Range("H90").Value = 6466.5835
$ws.Range("I90").Value = 2499.6667
$ws.Range("J90").Value = 7788.8887
$ws.Range("K90").Value = 22497.0003
$ws.Range("L90").Value = 70099.99830000001
$ws.Range("M90").Value = -16257.0003
$ws.Range("N90").Value = -82579.99830000001

$ws.Range("H107").Value = 1063.6786
$ws.Range("J107").Value = 1623.75
$ws.Range("L107").Value = 4871.25
$ws.Range("N107").Value = -8711.25

$ws.Range("H113").Value = 655.6949
$ws.Range("I113").Value = 694.5946
$ws.Range("J113").Value = 590.2727
$ws.Range("K113").Value = 2083.7838
$ws.Range("L113").Value = 1770.8181
$ws.Range("M113").Value = 86.21619999999984
$ws.Range("N113").Value = -6110.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7500
$ws.Range("J80").Value = 7500
$ws.Range("L80").Value = 7500
$ws.Range("N80").Value = -9496

$ws.Range("H83").Value = 7500
$ws.Range("J83").Value = 7500
$ws.Range("L83").Value = 37500
$ws.Range("N83").Value = -47484

$ws.Range("H122").Value = 5960.5835
$ws.Range("I122").Value = 7714.875
$ws.Range("J122").Value = 2452
$ws.Range("K122").Value = 23144.625
$ws.Range("L122").Value = 7356
$ws.Range("M122").Value = -20694.625
$ws.Range("N122").Value = -12256

$ws.Range("H132").Value = 2908.3809
$ws.Range("I132").Value = 3080.7693
$ws.Range("J132").Value = 2628.25
$ws.Range("K132").Value = 9242.3079
$ws.Range("L132").Value = 7884.75
$ws.Range("M132").Value = -6712.3079
$ws.Range("N132").Value = -12944.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5063
$ws.Range("I7").Value = 4321.1577
$ws.Range("J7").Value = 6824.875
$ws.Range("K7").Value = 4321.1577
$ws.Range("L7").Value = 6824.875
$ws.Range("M7").Value = -4209.1577
$ws.Range("N7").Value = -7048.875

$ws.Range("H40").Value = 3438.4119
$ws.Range("J40").Value = 3571.0833
$ws.Range("L40").Value = 3571.0833
$ws.Range("N40").Value = -3843.0833

$ws.Range("H75").Value = 36173
$ws.Range("J75").Value = 36173
$ws.Range("L75").Value = 36173
$ws.Range("N75").Value = -38045

$ws.Range("H78").Value = 36173
$ws.Range("J78").Value = 36173
$ws.Range("L78").Value = 108519
$ws.Range("N78").Value = -117879

$ws.Range("H126").Value = 5063
$ws.Range("I126").Value = 4321.1577
$ws.Range("J126").Value = 6824.875
$ws.Range("K126").Value = 12963.4731
$ws.Range("L126").Value = 20474.625
$ws.Range("M126").Value = -10493.4731
$ws.Range("N126").Value = -25414.625
